# Fill in the "Range" column (column 4) for several rows of the second
# table ("Type / Name / Display name / Range / Default / Required") with
# the documented value ranges for each parameter.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(2)

$ranges = @{
    2 = "1-100"    # studyName
    3 = "5-300"    # outPath
    4 = "16-128"   # access_key
    5 = "1-300"    # secret_key
    6 = "8-16"     # s3Region
    7 = "1-300"    # session_token
}

foreach ($rowIndex in $ranges.Keys) {
    $cell = $tbl.Cell($rowIndex, 4)
    $cell.Range.Text = $ranges[$rowIndex]
}
